$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Spp1"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 209.0063303333334
$ws.Range("H2").Value = 627.018991
$ws.Range("I2").Value = 0.6751081226665357
$ws.Range("J2").Value = 0.6751081226665357
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 1874.796815393856
$ws.Range("R2").Value = 16873.1713385447
$ws.Range("S2").Value = 0.3288510960319083
$ws.Range("T2").Value = 0.3288510960319082

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Spp1"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 209.0063303333334
$ws.Range("H3").Value = 627.018991
$ws.Range("I3").Value = 0.6751081226665357
$ws.Range("J3").Value = 0.6751081226665357
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 1883.579818744677
$ws.Range("R3").Value = 16952.21836870209
$ws.Range("S3").Value = 0.3303916897936715
$ws.Range("T3").Value = 0.3303916897936715

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Spp1"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 209.0063303333334
$ws.Range("H4").Value = 627.018991
$ws.Range("I4").Value = 0.6751081226665357
$ws.Range("J4").Value = 0.6751081226665357
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 90.4490918336159
$ws.Range("R4").Value = 814.041826502543
$ws.Range("S4").Value = 0.01586533684095609
$ws.Range("T4").Value = 0.01586533684095609

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Spp1"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.401741666666666
$ws.Range("H5").Value = 4.205225
$ws.Range("I5").Value = 0.004527744128790482
$ws.Range("J5").Value = 0.004527744128790482
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 12.5736900336
$ws.Range("R5").Value = 113.1632103024
$ws.Range("S5").Value = 0.002205503932353432
$ws.Range("T5").Value = 0.002205503932353432

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Spp1"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.401741666666666
$ws.Range("H6").Value = 4.205225
$ws.Range("I6").Value = 0.004527744128790482
$ws.Range("J6").Value = 0.004527744128790482
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 12.63259495641111
$ws.Range("R6").Value = 113.6933546077
$ws.Range("S6").Value = 0.002215836224508536
$ws.Range("T6").Value = 0.002215836224508536

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Spp1"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.401741666666666
$ws.Range("H7").Value = 4.205225
$ws.Range("I7").Value = 0.004527744128790482
$ws.Range("J7").Value = 0.004527744128790482
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 0.6066144529361109
$ws.Range("R7").Value = 5.459530076425
$ws.Range("S7").Value = 0.0001064039719285146
$ws.Range("T7").Value = 0.0001064039719285146

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Spp1"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 99.18134533333334
$ws.Range("H8").Value = 297.544036
$ws.Range("I8").Value = 0.3203641332046738
$ws.Range("J8").Value = 0.3203641332046737
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.970048
$ws.Range("N8").Value = 26.910144
$ws.Range("O8").Value = 0.487108783009476
$ws.Range("P8").Value = 0.4871087830094759
$ws.Range("Q8").Value = 889.6614283445761
$ws.Range("R8").Value = 8006.952855101184
$ws.Range("S8").Value = 0.1560521830452143
$ws.Range("T8").Value = 0.1560521830452143

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Spp1"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 99.18134533333334
$ws.Range("H9").Value = 297.544036
$ws.Range("I9").Value = 0.3203641332046738
$ws.Range("J9").Value = 0.3203641332046737
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.012070666666666
$ws.Range("N9").Value = 27.036212
$ws.Range("O9").Value = 0.489390778604016
$ws.Range("P9").Value = 0.489390778604016
$ws.Range("Q9").Value = 893.8292929590702
$ws.Range("R9").Value = 8044.463636631632
$ws.Range("S9").Value = 0.156783252585836
$ws.Range("T9").Value = 0.1567832525858359

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Spp1"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 99.18134533333334
$ws.Range("H10").Value = 297.544036
$ws.Range("I10").Value = 0.3203641332046738
$ws.Range("J10").Value = 0.3203641332046737
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4327576666666667
$ws.Range("N10").Value = 1.298273
$ws.Range("O10").Value = 0.02350043838650813
$ws.Range("P10").Value = 0.02350043838650813
$ws.Range("Q10").Value = 42.92148758331422
$ws.Range("R10").Value = 386.293388249828
$ws.Range("S10").Value = 0.007528697573623518
$ws.Range("T10").Value = 0.007528697573623517

